$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 13.29643343826742
$ws.Range("C2").Value = 8.430627994093014
$ws.Range("E2").Value = 23.97950684076044
$ws.Range("F2").Value = 39.71623611335701
$ws.Range("G2").Value = 3.608558167608966
$ws.Range("J2").Value = 7.627427706543941
$ws.Range("N2").Value = 16.16378226076835
$ws.Range("O2").Value = 19.3859581748509

$ws.Range("B3").Value = 12.71109282124262
$ws.Range("C3").Value = 7.965083392051193
$ws.Range("E3").Value = 23.73300734455565
$ws.Range("F3").Value = 39.46798755326568
$ws.Range("G3").Value = 3.610708912640057
$ws.Range("J3").Value = 7.652217005600419
$ws.Range("N3").Value = 16.20240832673113
$ws.Range("O3").Value = 19.4408127813896

$ws.Range("B4").Value = 12.33892020364333
$ws.Range("C4").Value = 7.663476833627128
$ws.Range("E4").Value = 23.58574254925725
$ws.Range("F4").Value = 39.32693451360839
$ws.Range("G4").Value = 3.612098313639028
$ws.Range("J4").Value = 7.668456600434549
$ws.Range("N4").Value = 16.22803490242455
$ws.Range("O4").Value = 19.48068904485964

$ws.Range("B5").Value = 12.18425567004727
$ws.Range("C5").Value = 7.53666437833051
$ws.Range("E5").Value = 23.52681848591964
$ws.Range("F5").Value = 39.2723621739228
$ws.Range("G5").Value = 3.612681872374225
$ws.Range("J5").Value = 7.675330809649552
$ws.Range("N5").Value = 16.23895887137366
$ws.Range("O5").Value = 19.49848854592505

$ws.Range("B6").Value = 12.15839919009058
$ws.Range("C6").Value = 7.515373344440763
$ws.Range("E6").Value = 23.51710157917839
$ws.Range("F6").Value = 39.26347741201953
$ws.Range("G6").Value = 3.612779822444425
$ws.Range("J6").Value = 7.676487762550795
$ws.Range("N6").Value = 16.24080185424786
$ws.Range("O6").Value = 19.50153748026621

$ws.Range("B7").Value = 12.3368461963558
$ws.Range("C7").Value = 7.661782323089255
$ws.Range("E7").Value = 23.58494339797698
$ws.Range("F7").Value = 39.32618669834637
$ws.Range("G7").Value = 3.612106113323874
$ws.Range("J7").Value = 7.668548269774631
$ws.Range("N7").Value = 16.22818027877197
$ws.Range("O7").Value = 19.48092283222775

$ws.Range("B8").Value = 13.09737604738587
$ws.Range("C8").Value = 8.273409960007434
$ws.Range("E8").Value = 23.89370214950011
$ws.Range("F8").Value = 39.62831016611347
$ws.Range("G8").Value = 3.609285492483409
$ws.Range("J8").Value = 7.635763735411075
$ws.Range("N8").Value = 16.17670448773974
$ws.Range("O8").Value = 19.40358138787666

$ws.Range("B9").Value = 14.47958729923316
$ws.Range("C9").Value = 9.345829300482752
$ws.Range("E9").Value = 24.52887814728147
$ws.Range("F9").Value = 40.30864149948159
$ws.Range("G9").Value = 3.604297810945738
$ws.Range("J9").Value = 7.579546805367922
$ws.Range("N9").Value = 16.09088787019984
$ws.Range("O9").Value = 19.30141869028296

$ws.Range("B10").Value = 15.41974693506985
$ws.Range("C10").Value = 10.05442691315731
$ws.Range("E10").Value = 25.0098915005126
$ws.Range("F10").Value = 40.85854236493695
$ws.Range("G10").Value = 3.600961039898994
$ws.Range("J10").Value = 7.543151183078006
$ws.Range("N10").Value = 16.03702343867038
$ws.Range("O10").Value = 19.25697737622461

$ws.Range("B11").Value = 15.8296622993498
$ws.Range("C11").Value = 10.35931789938339
$ws.Range("E11").Value = 25.23101521451716
$ws.Range("F11").Value = 41.11874131836648
$ws.Range("G11").Value = 3.599513423963897
$ws.Range("J11").Value = 7.527656382478338
$ws.Range("N11").Value = 16.01450601151682
$ws.Range("O11").Value = 19.24348952816879

$ws.Range("B12").Value = 15.98223668311627
$ws.Range("C12").Value = 10.47224873692964
$ws.Range("E12").Value = 25.31501200667959
$ws.Range("F12").Value = 41.21864097548888
$ws.Range("G12").Value = 3.598975298650953
$ws.Range("J12").Value = 7.521941387222516
$ws.Range("N12").Value = 16.00626423211462
$ws.Range("O12").Value = 19.23935520360026

$ws.Range("B13").Value = 15.94949636906726
$ws.Range("C13").Value = 10.44803953988268
$ws.Range("E13").Value = 25.29691124893018
$ws.Range("F13").Value = 41.19706620781484
$ws.Range("G13").Value = 3.599090747192636
$ws.Range("J13").Value = 7.523165430984043
$ws.Range("N13").Value = 16.00802657505002
$ws.Range("O13").Value = 19.24020223736849

$ws.Range("B14").Value = 15.8422683279259
$ws.Range("C14").Value = 10.36865945681052
$ws.Range("E14").Value = 25.2379208284607
$ws.Range("F14").Value = 41.12693315676579
$ws.Range("G14").Value = 3.599468950849488
$ws.Range("J14").Value = 7.527183149767003
$ws.Range("N14").Value = 16.01382224455709
$ws.Range("O14").Value = 19.24312986000565

$ws.Range("B15").Value = 15.77624012707055
$ws.Range("C15").Value = 10.31970769983746
$ws.Range("E15").Value = 25.2018195499086
$ws.Range("F15").Value = 41.08415046147086
$ws.Range("G15").Value = 3.599701919445464
$ws.Range("J15").Value = 7.529663981301198
$ws.Range("N15").Value = 16.01740937176126
$ws.Range("O15").Value = 19.24505001217305

$ws.Range("B16").Value = 15.39259184284641
$ws.Range("C16").Value = 10.03414897338152
$ws.Range("E16").Value = 24.99548101348374
$ws.Range("F16").Value = 40.84173342118708
$ws.Range("G16").Value = 3.601057055012237
$ws.Range("J16").Value = 7.544185158637171
$ws.Range("N16").Value = 16.03853492385204
$ws.Range("O16").Value = 19.25799473707745

$ws.Range("B17").Value = 15.15261140136323
$ws.Range("C17").Value = 9.854485495186822
$ws.Range("E17").Value = 24.86944041124466
$ws.Range("F17").Value = 40.69553791848438
$ws.Range("G17").Value = 3.601906353717625
$ws.Range("J17").Value = 7.553365275692404
$ws.Range("N17").Value = 16.05200303154051
$ws.Range("O17").Value = 19.26766353876945

$ws.Range("B18").Value = 15.01291329322171
$ws.Range("C18").Value = 9.749505559214853
$ws.Range("E18").Value = 24.79716663821511
$ws.Range("F18").Value = 40.61240053298995
$ws.Range("G18").Value = 3.602401467810138
$ws.Range("J18").Value = 7.558745380331352
$ws.Range("N18").Value = 16.05993646848459
$ws.Range("O18").Value = 19.27385782378728

$ws.Range("B19").Value = 14.96533068912063
$ws.Range("C19").Value = 9.713679630951733
$ws.Range("E19").Value = 24.7727361486367
$ws.Range("F19").Value = 40.58441716241008
$ws.Range("G19").Value = 3.602570243546779
$ws.Range("J19").Value = 7.560584160900374
$ws.Range("N19").Value = 16.06265471517795
$ws.Range("O19").Value = 19.27606364487234

$ws.Range("B20").Value = 15.17833108002334
$ws.Range("C20").Value = 9.873781069384746
$ws.Range("E20").Value = 24.88283523681453
$ws.Range("F20").Value = 40.71100286241276
$ws.Range("G20").Value = 3.601815259668173
$ws.Range("J20").Value = 7.552377693681477
$ws.Range("N20").Value = 16.05054998457154
$ws.Range("O20").Value = 19.26656871721316

$ws.Range("B21").Value = 15.87383646219848
$ws.Range("C21").Value = 10.39204392930752
$ws.Range("E21").Value = 25.25524118831507
$ws.Range("F21").Value = 41.14749643322082
$ws.Range("G21").Value = 3.599357590789628
$ws.Range("J21").Value = 7.525998909203792
$ws.Range("N21").Value = 16.01211218258872
$ws.Range("O21").Value = 19.24224349160928

$ws.Range("B22").Value = 16.31289912715328
$ws.Range("C22").Value = 10.71603853634981
$ws.Range("E22").Value = 25.50012531278899
$ws.Range("F22").Value = 41.44070625267886
$ws.Range("G22").Value = 3.597809948314777
$ws.Range("J22").Value = 7.509647970253694
$ws.Range("N22").Value = 15.98865232076035
$ws.Range("O22").Value = 19.23202039222123

$ws.Range("B23").Value = 16.08000851787396
$ws.Range("C23").Value = 10.54446720845308
$ws.Range("E23").Value = 25.36931219919362
$ws.Range("F23").Value = 41.2835145561618
$ws.Range("G23").Value = 3.598630611117458
$ws.Range("J23").Value = 7.51829345964293
$ws.Range("N23").Value = 16.00102141876752
$ws.Range("O23").Value = 19.23695572514422

$ws.Range("B24").Value = 15.16670859483023
$ws.Range("C24").Value = 9.865062800077361
$ws.Range("E24").Value = 24.87677884319436
$ws.Range("F24").Value = 40.70400831340456
$ws.Range("G24").Value = 3.601856421974101
$ws.Range("J24").Value = 7.552823860737108
$ws.Range("N24").Value = 16.05120631382895
$ws.Range("O24").Value = 19.26706170666549

$ws.Range("B25").Value = 14.1183465998811
$ws.Range("C25").Value = 9.069533522677675
$ws.Range("E25").Value = 24.35425001319811
$ws.Range("F25").Value = 40.11553111375434
$ws.Range("G25").Value = 3.605589302469629
$ws.Range("J25").Value = 7.593892223938979
$ws.Range("N25").Value = 16.11248793172099
$ws.Range("O25").Value = 19.32370902538016
